$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to text format so numeric-looking strings
# (e.g. "596.36") are preserved exactly and are not coerced into floating
# point numbers by the COM layer.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "72.224.90"
$ws.Range("E2").Value = "  -0.18%  "
$ws.Range("D3").Value = "2.650.89"
$ws.Range("E3").Value = "  +1.14%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "596.36"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").Value = "174.42"
$ws.Range("E6").Value = "  -2.29%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("D8").Value = "0.523"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").Value = "2.649.00"
$ws.Range("E9").Value = "  +1.09%  "
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("E11").Value = "  +2.18%  "
$ws.Range("D12").Value = "0.357"
$ws.Range("E12").Value = "  +0.97%  "
$ws.Range("E13").Value = "  -0.83%  "
$ws.Range("D14").Value = "3.140.36"
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("B15").Value = "WrappedBTC"
$ws.Range("C15").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D15").Value = "72.129.08"
$ws.Range("E15").Value = "  -0.22%  "
$ws.Range("B16").Value = "ShibaInu"
$ws.Range("C16").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D16").Value = "0.0000185"
$ws.Range("E16").Value = "  -1.81%  "
$ws.Range("D17").Value = "26.24"
$ws.Range("E17").Value = "  -1.49%  "
$ws.Range("D18").Value = "2.651.65"
$ws.Range("E18").Value = "  +0.94%  "
$ws.Range("D19").Value = "12.20"
$ws.Range("E19").Value = "  +5.14%  "
$ws.Range("D20").Value = "8.12"
$ws.Range("E20").Value = "  +1.47%  "
$ws.Range("D21").Value = "369.84"
$ws.Range("E21").Value = "  -2.68%  "
$ws.Range("D22").Value = "4.17"
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("D23").Value = "2.04"
$ws.Range("E23").Value = "  +0.34%  "
$ws.Range("D24").Value = "72.09"
$ws.Range("E24").Value = "  -1.72%  "
$ws.Range("D25").Value = "0.999"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("D26").Value = "4.31"
$ws.Range("E26").Value = "  -1.81%  "
$ws.Range("D27").Value = "9.73"
$ws.Range("E27").Value = "  -2.49%  "
$ws.Range("D28").Value = "2.792.95"
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("D30").Value = "0.0₃0967"
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("D31").Value = "8.08"
$ws.Range("E31").Value = "  -0.14%  "
$ws.Range("D32").Value = "498.46"
$ws.Range("E32").Value = "  -4.13%  "
$ws.Range("E33").Value = "  -2.71%  "
$ws.Range("D34").Value = "1.82"
$ws.Range("E34").Value = "  -0.60%  "
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").Value = "162.82"
$ws.Range("E36").Value = "  -1.63%  "
$ws.Range("D37").Value = "19.48"
$ws.Range("E37").Value = "  +0.85%  "
$ws.Range("B38").Value = "WhiteBITCoin"
$ws.Range("C38").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D38").Value = "18.92"
$ws.Range("E38").Value = "  -0.78%  "
$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").Value = "0.111"
$ws.Range("E39").Value = "  +0.40%  "
$ws.Range("E40").Value = "  -1.92%  "
$ws.Range("D41").Value = "1.76"
$ws.Range("E41").Value = "  -4.09%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("B43").Value = "RenderToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D43").Value = "4.98"
$ws.Range("E43").Value = "  -1.87%  "
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "2.58"
$ws.Range("E44").Value = "  -0.13%  "
$ws.Range("D45").Value = "0.331"
$ws.Range("E45").Value = "  -0.14%  "
$ws.Range("D46").Value = "39.45"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").Value = "155.61"
$ws.Range("E47").Value = "  +4.04%  "
$ws.Range("D48").Value = "3.73"
$ws.Range("E48").Value = "  +0.75%  "
$ws.Range("D49").Value = "0.556"
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("D50").Value = "1.72"
$ws.Range("E50").Value = "  +1.75%  "
$ws.Range("D51").Value = "0.0756"
$ws.Range("E51").Value = "  -1.27%  "

# Restore the default (unstyled) cell style now that the text values are set,
# matching the original workbook where these cells carry no explicit style.
$ws.Range("D2:D51").Style = "Normal"
